$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1 (rows 2-6): CNNBaseline ---
$ws.Range("B2").Value = "CNNBaseline"
$ws.Range("C2").Value = "107.83s"
$ws.Range("D2").Value = "4.26s"

$ws.Range("B3").Value = "CNNBaseline"
$ws.Range("C3").Value = "9.92s"
$ws.Range("D3").Value = "0.39s"

$ws.Range("B4").Value = "CNNBaseline"
$ws.Range("C4").Value = 0.85060000000000002
$ws.Range("D4").Value = 0.85440000000000005

$ws.Range("B5").Value = "5.75MB"

$ws.Range("A6").Value = "Parameters"
$ws.Range("B6").Value = 1507706

# --- Block 2 (rows 8-12): CNN Modified ---
$ws.Range("B8").Value = "CNN Modified"
$ws.Range("C8").Value = "192.54s"
$ws.Range("D8").Value = "6.21s"

$ws.Range("B9").Value = "CNN Modified"
$ws.Range("C9").Value = "16.09s"
$ws.Range("D9").Value = "0.49s"

$ws.Range("B10").Value = "CNN Modified"
$ws.Range("C10").Value = 0.87649999999999995
$ws.Range("D10").Value = 0.88090000000000002

$ws.Range("B11").Value = "3.13MB"

$ws.Range("A12").Value = "Parameters"
$ws.Range("B12").Value = 821130

# --- Block 3 (rows 14-18): Resnet18 modified ---
$ws.Range("B14").Value = "Resnet18 modified"
$ws.Range("C14").Value = "449.42s"
$ws.Range("D14").Value = "11.05s"

$ws.Range("B15").Value = "Resnet18 modified"
$ws.Range("C15").Value = "24.78s"
$ws.Range("D15").Value = "0.72s"

$ws.Range("B16").Value = "Resnet18 modified"
$ws.Range("C16").Value = 0.86
$ws.Range("D16").Value = 0.8619

$ws.Range("B17").Value = "42.65MB"

$ws.Range("A18").Value = "Parameters"
$ws.Range("B18").Value = 11181642

# --- Block 4 (rows 20-24): ViT (unchanged data, plus new Parameters row) ---
$ws.Range("A24").Value = "Parameters"
$ws.Range("B24").Value = 85806346

# Update selection to match target
$ws.Range("B7").Select()
